$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6039.1304
$ws.Range("J17").Value = 6245.4546
$ws.Range("L17").Value = 18736.3638
$ws.Range("N17").Value = -19072.3638

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8392.25
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8392.25
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8392.25
$ws.Range("N51").Value = -9360.25
$ws.Range("M51").ClearContents()

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8570.4
$ws.Range("J62").Value = 8950.625
$ws.Range("L62").Value = 8950.625
$ws.Range("N62").Value = -10198.625

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8570.4
$ws.Range("J65").Value = 8950.625
$ws.Range("L65").Value = 44753.125
$ws.Range("N65").Value = -50993.125

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4745.476
$ws.Range("I131").Value = 3366.3572
$ws.Range("K131").Value = 10099.0716
$ws.Range("M131").Value = -5059.071599999999

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 814.32434
$ws.Range("I132").Value = 695
$ws.Range("K132").Value = 2085
$ws.Range("M132").Value = 445

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2964.8118
$ws.Range("I138").Value = 2048.3333
$ws.Range("J138").Value = 3161.2
$ws.Range("K138").Value = 6144.999899999999
$ws.Range("L138").Value = 9483.599999999999
$ws.Range("M138").Value = -1004.999899999999
$ws.Range("N138").Value = -19763.6

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2356.3125
$ws.Range("I141").Value = 2356.3125
$ws.Range("K141").Value = 7068.9375
$ws.Range("M141").Value = -1888.9375

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3347
$ws.Range("I97").Value = 3347
$ws.Range("K97").Value = 3347
$ws.Range("M97").Value = -2851

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28241.195
$ws.Range("I31").Value = 1939.6957
$ws.Range("J31").Value = 61848.668
$ws.Range("K31").Value = 1939.6957
$ws.Range("L31").Value = 61848.668
$ws.Range("M31").Value = -1644.6957
$ws.Range("N31").Value = -62438.668

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 28241.195
$ws.Range("I34").Value = 1939.6957
$ws.Range("J34").Value = 61848.668
$ws.Range("K34").Value = 1939.6957
$ws.Range("L34").Value = 61848.668
$ws.Range("M34").Value = -1737.6957
$ws.Range("N34").Value = -62252.668

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3245.5293
$ws.Range("I132").Value = 2676.92
$ws.Range("J132").Value = 4825
$ws.Range("K132").Value = 8030.76
$ws.Range("L132").Value = 14475
$ws.Range("M132").Value = -5500.76
$ws.Range("N132").Value = -19535

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 54999.5
$ws.Range("I47").Value = 100000
$ws.Range("J47").Value = 9999
$ws.Range("K47").Value = 300000
$ws.Range("L47").Value = 29997
$ws.Range("M47").Value = -299569
$ws.Range("N47").Value = -30859

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 7871
$ws.Range("I63").Value = 2299.5
$ws.Range("J63").Value = 19014
$ws.Range("K63").Value = 6898.5
$ws.Range("L63").Value = 57042
$ws.Range("M63").Value = -6149.5
$ws.Range("N63").Value = -58540

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 7871
$ws.Range("I66").Value = 2299.5
$ws.Range("J66").Value = 19014
$ws.Range("K66").Value = 20695.5
$ws.Range("L66").Value = 171126
$ws.Range("M66").Value = -16951.5
$ws.Range("N66").Value = -178614

# CUL row 82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 748
$ws.Range("I82").Value = 748
$ws.Range("K82").Value = 2244
$ws.Range("M82").Value = -1838

# CUL row 85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 748
$ws.Range("I85").Value = 748
$ws.Range("K85").Value = 2244
$ws.Range("M85").Value = -840

# CUL row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 5000
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 15000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -10910
$ws.Range("N110").ClearContents()

# CUL row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 12912
$ws.Range("I119").Value = 1000
$ws.Range("J119").Value = 14103.2
$ws.Range("K119").Value = 3000
$ws.Range("L119").Value = 42309.60000000001
$ws.Range("M119").Value = 1838
$ws.Range("N119").Value = -51985.60000000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6482075
$ws.Range("J122").Value = 7520437
$ws.Range("L122").Value = 67683933
$ws.Range("N122").Value = -67688833

# CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 41669444
$ws.Range("I126").Value = 1772.6666
$ws.Range("K126").Value = 5317.9998
$ws.Range("M126").Value = -377.9997999999996

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 11500465
$ws.Range("I134").Value = 5233.952
$ws.Range("J134").Value = 41675450
$ws.Range("K134").Value = 15701.856
$ws.Range("L134").Value = 125026350
$ws.Range("M134").Value = -10631.856
$ws.Range("N134").Value = -125036490

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14699.5
$ws.Range("J70").Value = 15600.4
$ws.Range("L70").Value = 15600.4
$ws.Range("N70").Value = -16140.4

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 14699.5
$ws.Range("J73").Value = 15600.4
$ws.Range("L73").Value = 15600.4
$ws.Range("N73").Value = -17472.4

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 336129.06
$ws.Range("I80").Value = 419161.4
$ws.Range("K80").Value = 419161.4
$ws.Range("M80").Value = -418163.4

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 336129.06
$ws.Range("I83").Value = 419161.4
$ws.Range("K83").Value = 2095807
$ws.Range("M83").Value = -2090815

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1685.8889
$ws.Range("J97").Value = 2549.75
$ws.Range("L97").Value = 2549.75
$ws.Range("N97").Value = -3541.75

# LTW row 45
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 50000
$ws.Range("J45").Value = 50000
$ws.Range("L45").Value = 50000
$ws.Range("N45").Value = -50814

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7650.7036
$ws.Range("I46").Value = 7527.5
$ws.Range("J46").Value = 7685.905
$ws.Range("K46").Value = 7527.5
$ws.Range("L46").Value = 7685.905
$ws.Range("M46").Value = -7339.5
$ws.Range("N46").Value = -8061.905

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 10023
$ws.Range("J48").Value = 10023
$ws.Range("L48").Value = 10023
$ws.Range("N48").Value = -11345

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1472646
$ws.Range("I55").Value = 2632757.2
$ws.Range("K55").Value = 2632757.2
$ws.Range("M55").Value = -2632584.2

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3539
$ws.Range("J68").Value = 2998.6667
$ws.Range("L68").Value = 2998.6667
$ws.Range("N68").Value = -4496.6667

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3539
$ws.Range("J71").Value = 2998.6667
$ws.Range("L71").Value = 14993.3335
$ws.Range("N71").Value = -22481.3335

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6670.3
$ws.Range("J82").Value = 18618.666
$ws.Range("L82").Value = 18618.666
$ws.Range("N82").Value = -19340.666

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 6670.3
$ws.Range("J85").Value = 18618.666
$ws.Range("L85").Value = 18618.666
$ws.Range("N85").Value = -21114.666

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4722.3394
$ws.Range("I132").Value = 4753.5
$ws.Range("J132").Value = 4608.0835
$ws.Range("K132").Value = 14260.5
$ws.Range("L132").Value = 13824.2505
$ws.Range("M132").Value = -11730.5
$ws.Range("N132").Value = -18884.2505

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11143.286
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 11143.286
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1909.122
$ws.Range("I132").Value = 1450.4
$ws.Range("J132").Value = 4585
$ws.Range("K132").Value = 4351.200000000001
$ws.Range("L132").Value = 13755
$ws.Range("M132").Value = -1821.200000000001
$ws.Range("N132").Value = -18815
